# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.300.82'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').Value = '1.928.60'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7176'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.34%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.11'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.78%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3203'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.47%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07095'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7900'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.19%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07999'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.86%  '

$ws.Range('D13').Value = '1.930.81'
$ws.Range('E13').Value = '  -0.50%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.376'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.55%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.68%  '

$ws.Range('D17').Value = '30.300.51'
$ws.Range('E17').Value = '  -0.18%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '257.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.36%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008092'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.64%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.761'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.89%  '

$ws.Range('D21').Value = '2.182.23'
$ws.Range('E21').Value = '  -0.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.05%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.826'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.546'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.12%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.58'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.04%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.50%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.271'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.84%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1274'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.354'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.77%  '

$ws.Range('E31').Value = '  -2.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.402'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.89%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.142'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.49%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05125'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.77%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.273'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.01%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7464'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.85%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.774'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.10%  '

$ws.Range('E38').Value = '  -0.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.798'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.40%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.367'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.85%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4510'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8439'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.77%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9996'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.05%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.781'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.22%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.468'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.58%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.89%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '952.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.57%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4209'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.44%  '
